$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 44330
$ws.Range("D5").Value = 44316
$ws.Range("D6").Value = 44313
$ws.Range("M6").Value = 120
$ws.Range("Q6").Value = '$/caja 10 kilos empedrada'
$ws.Range("S6").Value = 11500
$ws.Range("T6").Value = 1
$ws.Range("D7").Value = 44306
$ws.Range("M7").Value = 80
$ws.Range("D8").Value = 44309
$ws.Range("Q8").Value = '$/caja 14 kilos granel'
$ws.Range("S8").Value = 821
$ws.Range("T8").Value = 14
$ws.Range("D9").Value = 44322
$ws.Range("M9").Value = 60
$ws.Range("D10").Value = 44302
